$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = "Biker mice"
$ws.Range("B25").Value = "Nicolas Giordani  | FC SAVIGNANO"
$ws.Range("C25").Value = "Federico Fasanelli | SBARX"
$ws.Range("D25").Value = "Matteo Diener | U.SGUARNA"
$ws.Range("E25").Value = "Filippo Benetti | I Magnifici"
$ws.Range("F25").Value = "Alessio Debiasi | Mai una gioia"
